$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: assign a text-formula, then copy/paste-special as values only so the
# resulting cell is a genuine text (string) value -- matching the original inlineStr cells --
# instead of letting Excel auto-coerce numeric-looking text ("327.34", "1.000", ...) into a
# floating point number.

$ws.Range("D2").Formula = '="30.135.72"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +5.79%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="1.920.78"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +2.48%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("E4").Formula = '="  -0.98%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="327.34"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +3.51%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("E6").Formula = '="  -0.94%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("D7").Formula = '="0.5174"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Formula = '="  +1.77%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

$ws.Range("D8").Formula = '="0.4049"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="  +4.07%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("D9").Formula = '="0.08466"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +1.11%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

$ws.Range("B10").Formula = '="Polygon"'
$ws.Range("B10").Copy()
$ws.Range("B10").PasteSpecial(-4163)
$ws.Range("C10").Formula = '="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"'
$ws.Range("C10").Copy()
$ws.Range("C10").PasteSpecial(-4163)
$ws.Range("D10").Formula = '="1.127"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +2.17%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$ws.Range("B11").Formula = '="OKB"'
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$ws.Range("C11").Formula = '="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"'
$ws.Range("C11").Copy()
$ws.Range("C11").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="42.81"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="  +2.49%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("D12").Formula = '="22.31"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  +9.33%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

$ws.Range("D13").Formula = '="6.356"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +2.23%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

$ws.Range("D14").Formula = '="1.925.84"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="  +2.94%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

$ws.Range("D15").Formula = '="7.377"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  +1.98%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("E16").Formula = '="  -1.09%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("D17").Formula = '="96.09"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +5.25%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

$ws.Range("D18").Formula = '="0.00001117"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = '="  +1.27%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("D19").Formula = '="0.06735"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  +0.05%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("D20").Formula = '="18.28"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  +3.25%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("E21").Formula = '="  -0.88%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

$ws.Range("D22").Formula = '="6.064"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  +2.28%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

$ws.Range("D23").Formula = '="30.128.91"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  +5.65%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

$ws.Range("D24").Formula = '="11.31"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = '="  +2.16%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("D25").Formula = '="2.201"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  -1.58%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

$ws.Range("D26").Formula = '="2.146.36"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +2.92%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

$ws.Range("B27").Formula = '="EthereumClassic"'
$ws.Range("B27").Copy()
$ws.Range("B27").PasteSpecial(-4163)
$ws.Range("C27").Formula = '="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"'
$ws.Range("C27").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="21.32"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  +3.42%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("B28").Formula = '="Monero"'
$ws.Range("B28").Copy()
$ws.Range("B28").PasteSpecial(-4163)
$ws.Range("C28").Formula = '="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"'
$ws.Range("C28").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D28").Formula = '="160.83"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = '="  -0.62%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$ws.Range("D29").Formula = '="2.463"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = '="  +3.53%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

$ws.Range("D30").Formula = '="129.26"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  +2.74%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("D31").Formula = '="1.085"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = '="  +4.65%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

$ws.Range("D32").Formula = '="0.1060"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = '="  +1.48%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

$ws.Range("D33").Formula = '="6.087"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = '="  +5.52%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

$ws.Range("D34").Formula = '="3.660"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +1.10%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

$ws.Range("D35").Formula = '="0.02516"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  +2.34%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$ws.Range("D36").Formula = '="0.06607"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  +0.81%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

$ws.Range("D37").Formula = '="0.2213"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  +2.41%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

$ws.Range("D38").Formula = '="1.239"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  +4.30%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

$ws.Range("D39").Formula = '="5.224"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  +3.08%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

$ws.Range("D40").Formula = '="9.033"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  +2.22%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

$ws.Range("D41").Formula = '="0.6578"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +2.81%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

$ws.Range("B42").Formula = '="Aptos"'
$ws.Range("B42").Copy()
$ws.Range("B42").PasteSpecial(-4163)
$ws.Range("C42").Formula = '="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"'
$ws.Range("C42").Copy()
$ws.Range("C42").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="11.76"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +5.93%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

$ws.Range("B43").Formula = '="TrustWalletToken"'
$ws.Range("B43").Copy()
$ws.Range("B43").PasteSpecial(-4163)
$ws.Range("C43").Formula = '="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"'
$ws.Range("C43").Copy()
$ws.Range("C43").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="1.251"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  -0.03%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

$ws.Range("D44").Formula = '="0.6177"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +2.76%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

$ws.Range("D45").Formula = '="13.18"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  +1.57%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

$ws.Range("E46").Formula = '="  +1.82%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

$ws.Range("D47").Formula = '="2.070"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +3.17%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("D48").Formula = '="1.246"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  +2.71%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

$ws.Range("D49").Formula = '="125.74"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = '="  +3.20%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

$ws.Range("D50").Formula = '="1.161"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +3.06%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

$ws.Range("D51").Formula = '="79.56"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  +4.25%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = 0
